# Swap the order of names in the "Recorded By" column (G) so that
# "dnasr281@gmail.com" is listed first whenever it appears together
# with exactly one other recorder (e.g. "System" or "admin@admin.com").
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#          "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -eq 2 -and ($parts[0] -eq $target -or $parts[1] -eq $target) -and $parts[0] -ne $parts[1]) {
            if ($parts[0] -ne $target) {
                $newVal = $parts[1] + ", " + $parts[0]
                $cell.Value2 = $newVal
            }
        }
    }
}

Write-Host "Recorded By swap complete"
